$excel = New-Object -ComObject Excel.Application
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = 100.1
$ws.Range("B2").Value = 110.1
$ws.Range("C2").Value = 210.1

$ws.Range("A3").Value = 200.2
$ws.Range("B3").Value = 120.2
$ws.Range("C3").Value = 220.2
